# "selecting data" edit: duplicate the OS-table's label/value columns
# (A:B) into a second copy in columns C:D, for every row that holds
# actual data (the section-header rows such as "OS1", "OS2", ... only
# have a label in column A and are left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 41; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    if ($bVal -ne $null) {
        $aVal = $ws.Cells.Item($r, 1).Value()
        $ws.Cells.Item($r, 3).Value = $aVal
        $ws.Cells.Item($r, 4).Value = $bVal
    }
}

# Match the author's saved view/selection state: window scrolled down to
# row 19, with the new C38:D41 block selected (C38 active).
$ws.Range("C38:D41").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
